$d = $word.ActiveDocument

# 1. Fix the typo "net" -> "new" in the "Add a link to the help website..." bullet.
$old = "Add a link to the help website to the spreadsheet. This could be a net sheet named help. The help sheet should be the first sheet that is shown when the template is first opened."
$new = "Add a link to the help website to the spreadsheet. This could be a new sheet named help. The help sheet should be the first sheet that is shown when the template is first opened."

$rng = $d.Content
$rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# 2. Move the (hidden) "_GoBack" bookmark from after "case." to right after
#    "This could be a new" (before the space that precedes "sheet named help").
#    Adding a bookmark with a name that already exists relocates it, which is
#    exactly the effect we want (and also splits the run at that point, just
#    like the authoring Word client did).
$marker = "Add a link to the help website to the spreadsheet. This could be a new"
$bmRange = $d.Content
$bmRange.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($bmRange.End, $bmRange.End)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
